# Update NATMI TPM-derived values for Efna4-Epha2 sheet.
# Both Efna4 (ligand) and Epha2 (receptor) expression in the "ECs" cluster
# changed, which cascades into the derived specificity/weight columns for
# every row of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E2" = 2
    "F2" = 0.6666666666666666
    "G2" = 0.2988413333333333
    "H2" = 0.896524
    "I2" = 0.3632971504731247
    "J2" = 0.3632971504731246
    "M2" = 10.858287
    "N2" = 32.574861
    "O2" = 0.5084025289165609
    "P2" = 0.508402528916561
    "Q2" = 3.244904964796
    "R2" = 29.204144683164
    "S2" = 0.1847011900487169
    "T2" = 0.184701190048717

    "E3" = 2
    "F3" = 0.6666666666666666
    "G3" = 0.2988413333333333
    "H3" = 0.896524
    "I3" = 0.3632971504731247
    "J3" = 0.3632971504731246
    "O3" = 0.004437346842596906
    "P3" = 0.004437346842596906
    "Q3" = 0.02832159161511111
    "R3" = 0.254894324536
    "S3" = 0.001612075463576373
    "T3" = 0.001612075463576373

    "E4" = 2
    "F4" = 0.6666666666666666
    "G4" = 0.2988413333333333
    "H4" = 0.896524
    "I4" = 0.3632971504731247
    "J4" = 0.3632971504731246
    "O4" = 0.4871601242408422
    "P4" = 0.4871601242408422
    "Q4" = 3.109324237958667
    "R4" = 27.983918141628
    "S4" = 0.1769838849608313
    "T4" = 0.1769838849608313

    "I5" = 0.4682720202225272
    "J5" = 0.4682720202225272
    "M5" = 10.858287
    "N5" = 32.574861
    "O5" = 0.5084025289165609
    "P5" = 0.508402528916561
    "Q5" = 4.182521666674999
    "R5" = 37.642695000075
    "S5" = 0.2380706793019998
    "T5" = 0.2380706793019998

    "I6" = 0.4682720202225272
    "J6" = 0.4682720202225272
    "O6" = 0.004437346842596906
    "P6" = 0.004437346842596906
    "S6" = 0.002077885370410906
    "T6" = 0.002077885370410906

    "I7" = 0.4682720202225272
    "J7" = 0.4682720202225272
    "O7" = 0.4871601242408422
    "P7" = 0.4871601242408422
    "S7" = 0.2281234555501165
    "T7" = 0.2281234555501165

    "I8" = 0.1684308293043481
    "J8" = 0.1684308293043481
    "M8" = 10.858287
    "N8" = 32.574861
    "O8" = 0.5084025289165609
    "P8" = 0.508402528916561
    "Q8" = 1.504393947276
    "R8" = 13.539545525484
    "S8" = 0.08563065956584419
    "T8" = 0.0856306595658442

    "I9" = 0.1684308293043481
    "J9" = 0.1684308293043481
    "O9" = 0.004437346842596906
    "P9" = 0.004437346842596906
    "S9" = 0.0007473860086096277
    "T9" = 0.0007473860086096277

    "I10" = 0.1684308293043481
    "J10" = 0.1684308293043481
    "O10" = 0.4871601242408422
    "P10" = 0.4871601242408422
    "S10" = 0.08205278372989432
    "T10" = 0.08205278372989433
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
